$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-9.858623504638672, -13.49863147735596, -4.903345584869385),
    @(1.233975768089294, 17.44783401489258, -5.1009521484375),
    @(6.964565753936768, -6.484397888183594, -1.327785611152649),
    @(-2.914164066314697, -17.45715522766113, -3.230347871780396),
    @(-0.99935120344162, 3.104047536849976, 2.052831172943115),
    @(-10.20110511779785, -13.58332061767578, -4.393616676330566),
    @(3.081676959991455, 6.77015495300293, 9.200222969055176),
    @(4.065448760986328, -7.618904590606689, -0.8031428456306458),
    @(-1.262471556663513, -17.20468711853027, -2.038849592208862),
    @(-0.1183775141835212, 2.852112531661988, 1.839245676994324)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
